# Updated cryptos list — refresh price / volume(1h) figures pulled from
# coinranking.com, and re-sync two pairs of rows whose underlying rank
# order flipped (Polkadot/Avalanche at rows 16-17, Maker/FraxShare/Aave
# at rows 46-48) so each row's Coin/Link/Price/Volume stay together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One entry per data row that changed: the column letters present in
# "Cells" are exactly the ones this refresh touched for that row.
$updates = @(
    @{ Row = 2; Cells = @{ D='36.353.30'; E='  +0.28%  ' } },
    @{ Row = 3; Cells = @{ D='2.018.11'; E='  -1.23%  ' } },
    @{ Row = 4; Cells = @{ E='  +0.06%  ' } },
    @{ Row = 5; Cells = @{ D='251.98'; E='  +2.99%  ' } },
    @{ Row = 6; Cells = @{ D='0.643'; E='  -2.64%  ' } },
    @{ Row = 7; Cells = @{ D='62.53'; E='  +10.72%  ' } },
    @{ Row = 8; Cells = @{ E='  -0.06%  ' } },
    @{ Row = 9; Cells = @{ D='59.35'; E='  -6.34%  ' } },
    @{ Row = 10; Cells = @{ D='0.371'; E='  +1.55%  ' } },
    @{ Row = 11; Cells = @{ D='0.0745'; E='  +0.17%  ' } },
    @{ Row = 12; Cells = @{ D='0.104'; E='  -1.50%  ' } },
    @{ Row = 13; Cells = @{ D='0.911'; E='  +0.13%  ' } },
    @{ Row = 14; Cells = @{ D='14.85'; E='  +5.21%  ' } },
    @{ Row = 15; Cells = @{ D='2.312.12'; E='  -1.23%  ' } },
    @{ Row = 16; Cells = @{ B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='20.38'; E='  +17.04%  ' } },
    @{ Row = 17; Cells = @{ B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='5.43'; E='  +1.23%  ' } },
    @{ Row = 18; Cells = @{ D='2.012.94'; E='  -1.20%  ' } },
    @{ Row = 19; Cells = @{ D='36.305.44'; E='  +0.21%  ' } },
    @{ Row = 20; Cells = @{ D='72.09'; E='  +1.20%  ' } },
    @{ Row = 21; Cells = @{ E='  +0.89%  ' } },
    @{ Row = 22; Cells = @{ D='5.31'; E='  +2.77%  ' } },
    @{ Row = 23; Cells = @{ D='234.57'; E='  -0.97%  ' } },
    @{ Row = 24; Cells = @{ D='2.65'; E='  +17.80%  ' } },
    @{ Row = 25; Cells = @{ E='  +0.12%  ' } },
    @{ Row = 26; Cells = @{ E='  -1.70%  ' } },
    @{ Row = 27; Cells = @{ D='9.63'; E='  +3.90%  ' } },
    @{ Row = 28; Cells = @{ D='163.24'; E='  -0.76%  ' } },
    @{ Row = 29; Cells = @{ D='19.66'; E='  -1.20%  ' } },
    @{ Row = 30; Cells = @{ E='  +33.69%  ' } },
    @{ Row = 31; Cells = @{ E='  -0.36%  ' } },
    @{ Row = 32; Cells = @{ D='5.13'; E='  +3.46%  ' } },
    @{ Row = 33; Cells = @{ E='  -1.06%  ' } },
    @{ Row = 34; Cells = @{ D='0.0610'; E='  +2.23%  ' } },
    @{ Row = 35; Cells = @{ D='4.54'; E='  +3.39%  ' } },
    @{ Row = 36; Cells = @{ E='  +12.57%  ' } },
    @{ Row = 37; Cells = @{ E='  +0.01%  ' } },
    @{ Row = 38; Cells = @{ E='  -0.36%  ' } },
    @{ Row = 39; Cells = @{ D='5.94'; E='  +17.88%  ' } },
    @{ Row = 40; Cells = @{ D='0.104'; E='  +14.60%  ' } },
    @{ Row = 41; Cells = @{ E='  +0.79%  ' } },
    @{ Row = 42; Cells = @{ D='2.92'; E='  +1.34%  ' } },
    @{ Row = 43; Cells = @{ E='  +0.80%  ' } },
    @{ Row = 44; Cells = @{ E='  +2.99%  ' } },
    @{ Row = 45; Cells = @{ D='16.63'; E='  +4.89%  ' } },
    @{ Row = 46; Cells = @{ B='Maker'; C='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D='1.439.40'; E='  +5.62%  ' } },
    @{ Row = 47; Cells = @{ B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='7.87'; E='  +6.83%  ' } },
    @{ Row = 48; Cells = @{ B='Aave'; C='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D='94.54'; E='  +1.48%  ' } },
    @{ Row = 49; Cells = @{ D='2.63'; E='  +15.59%  ' } },
    @{ Row = 50; Cells = @{ D='2.93'; E='  -0.59%  ' } },
    @{ Row = 51; Cells = @{ D='47.49'; E='  +4.17%  ' } }
)

foreach ($update in $updates) {
    $row = $update.Row
    foreach ($col in $update.Cells.Keys) {
        $cellRef = "$col$row"
        $value = $update.Cells[$col]
        $rng = $ws.Range($cellRef)

        # Force text storage (matches the source data's inlineStr type) even
        # for values that look numeric ("251.98", "20.38", ...) — plain
        # assignment would otherwise let Excel coerce those into Number
        # cells. Setting a "@" text format before the write, then clearing
        # formats afterward, keeps the cell's style back at its original
        # (default) index instead of leaving a stray number-format override.
        $rng.NumberFormat = "@"
        $rng.Value = $value
        $rng.ClearFormats()
    }
}
